# wdkCache fix: drop "+checksum" from the step-id box, rename
# user_dataset_id -> dataset_id, rename the old dataset_id box to
# "Content checksum", and retitle the deck "Params - Old".
#
# EMU <-> point helper. Shape.Width/Height/Left/Top are expressed in
# points (1 pt = 12700 EMU) and the interop stores them as single
# precision floats, so converting target-EMU/12700.0 directly rounds
# down by ~1 EMU once re-quantised on save. Nudging by half an EMU
# before the conversion lands exactly on the desired integer EMU
# value once PowerPoint re-serialises the xfrm.
function EMU-ToPoints($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1: title "Params - Current" -> "Params - Old"
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "Title 64") {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf(" - Current")
        if ($idx -ge 0) {
            $target = $tr.Characters($idx + 1, 10)
            $target.Text = " "
            [void]$tr.InsertAfter("- Old")
        }
    }
}

# ---------------------------------------------------------------
# Slide 2: wdk cache diagram boxes
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    if (-not $sh.TextFrame.HasText) { continue }

    $tr = $sh.TextFrame.TextRange
    $full = $tr.Text

    if ($sh.Name -eq "TextBox 32" -and $full -eq "step id+checksum") {
        # drop the "+checksum" suffix, keep "step id"
        $idx = $full.IndexOf("+checksum")
        $target = $tr.Characters($idx + 1, 9)
        $target.Text = ""
        $sh.Width = EMU-ToPoints(644728)
    }
    elseif ($sh.Name -eq "TextBox 57" -and $full -eq "user_dataset_id") {
        # drop the "user_" prefix, keep "dataset_id"
        $target = $tr.Characters(1, 5)
        $target.Text = ""
        $sh.Width = EMU-ToPoints(888385)
    }
    elseif ($sh.Name -eq "TextBox 62" -and $full -eq "dataset_id") {
        # replace the whole label with "Content checksum"
        $first = $tr.Characters(1, 1)
        $first.Text = ""
        $tr.Text = "Content checksum"
        $sh.Width = EMU-ToPoints(1449436)
    }
}
